$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that sometimes look numeric (e.g. "1.006",
# "0.3490"). Force those specific cells to Text format immediately before
# writing them so Excel does not silently coerce them into numbers (which
# would also drop significant trailing zeros / switch to scientific
# notation). Cells whose new text is not parseable as a plain number (e.g.
# "28.107.81", which has two dots) do not need the guard and are left alone
# so no stray formatting changes leak onto cells outside this edit.

$ws.Range("D2").Value = '28.107.81'
$ws.Range("E2").Value = '  +2.81%  '
$ws.Range("D3").Value = '1.804.59'
$ws.Range("E3").Value = '  +0.92%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.99'
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3933'
$ws.Range("E7").Value = '  +3.53%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3490'
$ws.Range("E8").Value = '  +1.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.25'
$ws.Range("E9").Value = '  -0.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.166'
$ws.Range("E10").Value = '  -2.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07529'
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.98'
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.503'
$ws.Range("E14").Value = '  +0.69%  '
$ws.Range("D15").Value = '1.806.66'
$ws.Range("E15").Value = '  +0.97%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.136'
$ws.Range("E16").Value = '  +0.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001099'
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06703'
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.90'
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.69'
$ws.Range("E21").Value = '  +1.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.551'
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").Value = '28.142.12'
$ws.Range("E23").Value = '  +2.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.39'
$ws.Range("E24").Value = '  -1.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.406'
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.482'
$ws.Range("E26").Value = '  -1.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.42'
$ws.Range("E27").Value = '  -0.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.513'
$ws.Range("E28").Value = '  -1.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.26'
$ws.Range("E29").Value = '  +0.56%  '
$ws.Range("D30").Value = '2.013.46'
$ws.Range("E30").Value = '  +0.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '135.22'
$ws.Range("E31").Value = '  +1.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.165'
$ws.Range("E32").Value = '  +1.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.022'
$ws.Range("E33").Value = '  -0.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08839'
$ws.Range("E34").Value = '  +1.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.07'
$ws.Range("E35").Value = '  -0.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6934'
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02422'
$ws.Range("E37").Value = '  +3.41%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06521'
$ws.Range("E38").Value = '  +2.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.432'
$ws.Range("E39").Value = '  -0.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.609'
$ws.Range("E40").Value = '  -2.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2210'
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("E42").Value = '  -0.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.465'
$ws.Range("E43").Value = '  -4.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.65'
$ws.Range("E44").Value = '  +1.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6400'
$ws.Range("E46").Value = '  -0.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.875'
$ws.Range("E47").Value = '  +0.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.142'
$ws.Range("E48").Value = '  +0.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.71'
$ws.Range("E49").Value = '  +0.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07189'
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.95'
$ws.Range("E51").Value = '  +0.50%  '
